$d = $word.ActiveDocument

# Locate the target paragraph: "Com base nos testes realizados ..."
$paras = $d.Paragraphs
$count = $paras.Count
$targetIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $paras.Item($i)
    $t = $p.Range.Text
    if ($t.Length -ge 10 -and $t.Substring(0, 10) -eq "Com base n") {
        $targetIndex = $i
    }
}

if ($targetIndex -eq -1) {
    Write-Output "Target paragraph not found"
} else {
    $p = $paras.Item($targetIndex)
    $pStart = $p.Range.Start

    # Insert "quase " right before "todas as funcionalidades ..."
    $insPos = $pStart + 72
    $insRange = $d.Range($insPos, $insPos)
    $insRange.InsertBefore("quase ")

    # Append a trailing space at the very end of the sentence (after "trabalho.")
    $p2 = $paras.Item($targetIndex)
    $p2.Range.InsertAfter(" ")

    # Insert a brand-new paragraph right after, matching the "PargrafodaLista"
    # style / indent / justification used elsewhere in this section.
    $p3 = $paras.Item($targetIndex)
    $p3.Range.InsertParagraphAfter()

    $paras2 = $d.Paragraphs
    $newPara = $paras2.Item($targetIndex + 1)
    $newPara.Range.InsertBefore("O item 6 cumpre a funcionalidade de printar o ranking de despesas por categoria em ordem decrescente de valor, mas após printar aparece uma mensagem de erro dizendo “adress out of range”, e após analisarmos entendemos que o erro está no syscall, mas não sabemos como resolve-lo.")

    Write-Output "edit applied"
}
